$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# New row 35: 0xaff221f9 / APP-Initer发现一个PPK模块 / 提示 / 内核-App
$ws.Range("A35").Value = "0xaff221f9"
$ws.Range("B35").Value = "APP-Initer发现一个PPK模块"
$ws.Range("C35").Value = "提示"
$ws.Range("D35").Value = "内核-App"

# New row 36: 0xaff221fa (trailing space) / APP-Initer发现一个EXT模块 / 提示 / 内核-App
$ws.Range("A36").Value = "0xaff221fa "
$ws.Range("B36").Value = "APP-Initer发现一个EXT模块"
$ws.Range("C36").Value = "提示"
$ws.Range("D36").Value = "内核-App"

# Update the visible selection to match the new bottom rows
$ws.Range("C36:D36").Select()
